$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF219").Value = 632.926916520989
$ws.Range("J220").Value = 28953.2297859252
$ws.Range("AF220").Value = 637.018899068649
$ws.Range("AF221").Value = 639.346200486882
$ws.Range("AF222").Value = 641.682004528663
$ws.Range("AF223").Value = 644.026342257697
$ws.Range("AF224").Value = 646.379244851179
$ws.Range("AF225").Value = 650.109646128924
$ws.Range("J226").Value = 30700.9838575899
$ws.Range("AF226").Value = 653.861576398825
$ws.Range("D227").Value = 0.00467044581501486
$ws.Range("E227").Value = 0.00620433826657418
$ws.Range("F227").Value = 0.00620433826657418
$ws.Range("G227").Value = 0.00620433826657418
$ws.Range("H227").Value = 0.00484384269748661
$ws.Range("D228").Value = 0.00493841965187825
$ws.Range("E228").Value = 0.00630060376804042
$ws.Range("F228").Value = 0.00630060376804042
$ws.Range("G228").Value = 0.00630060376804042
$ws.Range("H228").Value = 0.00479383672055
$ws.Range("D229").Value = 0.00504079228887644
$ws.Range("E229").Value = 0.00628133760843652
$ws.Range("F229").Value = 0.00628133760843652
$ws.Range("G229").Value = 0.00628133760843652
$ws.Range("H229").Value = 0.00474425523931421
$ws.Range("D230").Value = 0.00515873503879316
$ws.Range("E230").Value = 0.00630531064768802
$ws.Range("F230").Value = 0.00630531064768802
$ws.Range("G230").Value = 0.00630531064768802
$ws.Range("H230").Value = 0.00473037874429894
$ws.Range("D231").Value = 0.00511809259973983
$ws.Range("E231").Value = 0.00636942692478693
$ws.Range("F231").Value = 0.00636942692478693
$ws.Range("G231").Value = 0.00636942692478693
$ws.Range("H231").Value = 0.00473052086148029
$ws.Range("D232").Value = 0.00518873101404216
$ws.Range("E232").Value = 0.00641001309176392
$ws.Range("F232").Value = 0.00641001309176392
$ws.Range("G232").Value = 0.00641001309176392
$ws.Range("H232").Value = 0.00472642424019742
$ws.Range("D233").Value = 0.0055287138727893
$ws.Range("E233").Value = 0.0064507136259504
$ws.Range("F233").Value = 0.0064507136259504
$ws.Range("G233").Value = 0.0064507136259504
$ws.Range("H233").Value = 0.0047374828288711
$ws.Range("D234").Value = 0.0053958958468916
$ws.Range("E234").Value = 0.0065315402568018
$ws.Range("F234").Value = 0.0065315402568018
$ws.Range("G234").Value = 0.0065315402568018
$ws.Range("H234").Value = 0.004760735950051
$ws.Range("D235").Value = 0.005547011459617
$ws.Range("E235").Value = 0.00666145596157186
$ws.Range("F235").Value = 0.00666145596157186
$ws.Range("G235").Value = 0.00666145596157186
$ws.Range("H235").Value = 0.00479021974507288
$ws.Range("D236").Value = 0.00549985398617747
$ws.Range("E236").Value = 0.00674235738498696
$ws.Range("F236").Value = 0.00674235738498696
$ws.Range("G236").Value = 0.00674235738498696
$ws.Range("H236").Value = 0.004806028646549
$ws.Range("D237").Value = 0.00551817609758465
$ws.Range("E237").Value = 0.00683764486104899
$ws.Range("F237").Value = 0.00683764486104899
$ws.Range("G237").Value = 0.00683764486104899
$ws.Range("H237").Value = 0.00481903473851841
$ws.Range("D238").Value = 0.00547201687236676
$ws.Range("E238").Value = 0.00688738078154305
$ws.Range("F238").Value = 0.00688738078154305
$ws.Range("G238").Value = 0.00688738078154305
$ws.Range("H238").Value = 0.00483640202891866
$ws.Range("D239").Value = 0.00532843972065278
$ws.Range("E239").Value = 0.00688831459364736
$ws.Range("F239").Value = 0.00688831459364736
$ws.Range("G239").Value = 0.00688831459364736
$ws.Range("H239").Value = 0.00484025736739846
$ws.Range("D240").Value = 0.00539044209335926
$ws.Range("E240").Value = 0.00695124959344939
$ws.Range("F240").Value = 0.00695124959344939
$ws.Range("G240").Value = 0.00695124959344939
$ws.Range("H240").Value = 0.00484361940998301
$ws.Range("D241").Value = 0.0053714494404371
$ws.Range("E241").Value = 0.00698929263417969
$ws.Range("F241").Value = 0.00698929263417969
$ws.Range("G241").Value = 0.00698929263417969
$ws.Range("H241").Value = 0.00485075146532865
$ws.Range("D242").Value = 0.00516706709151271
$ws.Range("E242").Value = 0.00701085477153973
$ws.Range("F242").Value = 0.00701085477153973
$ws.Range("G242").Value = 0.00701085477153973
$ws.Range("H242").Value = 0.00485153149474238
$ws.Range("D243").Value = 0.0052944419865304
$ws.Range("E243").Value = 0.00705708111985137
$ws.Range("F243").Value = 0.00705708111985137
$ws.Range("G243").Value = 0.00705708111985137
$ws.Range("H243").Value = 0.0048632081837976
$ws.Range("D244").Value = 0.00526173532683871
$ws.Range("E244").Value = 0.00703501920168259
$ws.Range("F244").Value = 0.00703501920168259
$ws.Range("G244").Value = 0.00703501920168259
$ws.Range("H244").Value = 0.0048529979157681
$ws.Range("D245").Value = 0.00527563961968203
$ws.Range("E245").Value = 0.00706444317107557
$ws.Range("F245").Value = 0.00706444317107557
$ws.Range("G245").Value = 0.00706444317107557
$ws.Range("H245").Value = 0.00486232207935799
$ws.Range("D246").Value = 0.00530215486890606
$ws.Range("E246").Value = 0.00705152599487957
$ws.Range("F246").Value = 0.00705152599487957
$ws.Range("G246").Value = 0.00705152599487957
$ws.Range("H246").Value = 0.00485777293438971
$ws.Range("D247").Value = 0.00529234893313224
$ws.Range("E247").Value = 0.00711732623587236
$ws.Range("F247").Value = 0.00711732623587236
$ws.Range("G247").Value = 0.00711732623587236
$ws.Range("H247").Value = 0.00486115755082905
$ws.Range("D248").Value = 0.0052562450320528
$ws.Range("E248").Value = 0.00708977135022493
$ws.Range("F248").Value = 0.00708977135022493
$ws.Range("G248").Value = 0.00708977135022493
$ws.Range("H248").Value = 0.00486397246548265
$ws.Range("D249").Value = 0.00522846614789252
$ws.Range("E249").Value = 0.0071360334951307
$ws.Range("F249").Value = 0.0071360334951307
$ws.Range("G249").Value = 0.0071360334951307
$ws.Range("H249").Value = 0.00485809684592908
$ws.Range("D250").Value = 0.00493901013144238
$ws.Range("E250").Value = 0.00716903453216777
$ws.Range("F250").Value = 0.00716903453216777
$ws.Range("G250").Value = 0.00716903453216777
$ws.Range("H250").Value = 0.00486627672853213
$ws.Range("D251").Value = 0.0050454408249363
$ws.Range("E251").Value = 0.00718936727518638
$ws.Range("F251").Value = 0.00718936727518638
$ws.Range("G251").Value = 0.00718936727518638
$ws.Range("H251").Value = 0.00486546004897104
$ws.Range("J251").Value = 38588.5149279641
$ws.Range("D252").Value = 0.00508921192760492
$ws.Range("E252").Value = 0.00718313707806573
$ws.Range("F252").Value = 0.00718313707806573
$ws.Range("G252").Value = 0.00718313707806573
$ws.Range("H252").Value = 0.00486341439779814
$ws.Range("D253").Value = 0.00528501550993887
$ws.Range("E253").Value = 0.00717568925283363
$ws.Range("F253").Value = 0.00717568925283363
$ws.Range("G253").Value = 0.00717568925283363
$ws.Range("H253").Value = 0.00486299660338019
$ws.Range("D254").Value = 0.00530615579081561
$ws.Range("E254").Value = 0.00723850947089333
$ws.Range("F254").Value = 0.00723850947089333
$ws.Range("G254").Value = 0.00723850947089333
$ws.Range("H254").Value = 0.00486472626256185
$ws.Range("D255").Value = 0.00526366638494391
$ws.Range("E255").Value = 0.00720514060374988
$ws.Range("F255").Value = 0.00720514060374988
$ws.Range("G255").Value = 0.00720514060374988
$ws.Range("H255").Value = 0.00486977895366336
$ws.Range("D256").Value = 0.00533608914008776
$ws.Range("E256").Value = 0.00725110506577331
$ws.Range("F256").Value = 0.00725110506577331
$ws.Range("G256").Value = 0.00725110506577331
$ws.Range("H256").Value = 0.00486179575586099
$ws.Range("J256").Value = 40425.7205861868
$ws.Range("D257").Value = 0.00527908676509004
$ws.Range("E257").Value = 0.00729628986635267
$ws.Range("F257").Value = 0.00729628986635267
$ws.Range("G257").Value = 0.00729628986635267
$ws.Range("H257").Value = 0.00486832168670137
$ws.Range("D258").Value = 0.00528755943207071
$ws.Range("E258").Value = 0.00732436171924289
$ws.Range("F258").Value = 0.00732436171924289
$ws.Range("G258").Value = 0.00732436171924289
$ws.Range("H258").Value = 0.00486529484724985
$ws.Range("D259").Value = 0.00533846296484497
$ws.Range("E259").Value = 0.00736773064019669
$ws.Range("F259").Value = 0.00736773064019669
$ws.Range("G259").Value = 0.00736773064019669
$ws.Range("H259").Value = 0.00486656389504803
$ws.Range("D260").Value = 0.00530320809440532
$ws.Range("E260").Value = 0.00741037964071811
$ws.Range("F260").Value = 0.00741037964071811
$ws.Range("G260").Value = 0.00741037964071811
$ws.Range("H260").Value = 0.00486575776518472

Write-Output "Applied 182 cell updates"
